# DailyStatusTracker.xlsx update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("OpsTracker")
$ws.Activate()

# --- Update row 11 (Item 14): "College doors Mock test Dry run" -> "Initial College doors Mock test Dry run"
#     and add new comment to column E
$ws.Range("B11").Value = "Initial College doors Mock test Dry run"
$ws.Range("E11").Value = "Debashish and Ananya tested from our end and raised 25 defects. We notified College Doors"

# --- Update row 21 (Item 27) comment to have a trailing space
$ws.Range("E21").Value = "Why cannot teacher's type the question papers? "

# --- Insert a new row after row 21 (before old row 22) for a follow-up comment with no item number/description
$ws.Rows.Item(22).Insert()
$ws.Range("C22").Value = "Debasish"
$ws.Range("D22").Value = "WIP"
$ws.Range("E22").Value = "We communicated to Subroto Sir. He will do the necessary stuffs on Thursday (23-11-2023)"

# --- Append new row 30 at the bottom of the table
$ws.Range("A30").Value = 36
$ws.Range("B30").Value = "Dry run of Mock Test with 3 users "
$ws.Range("C30").Value = "Debasish"
$ws.Range("D30").Value = "Todo"
$ws.Range("E30").Value = "On 20th November 2023, we will test with 3 users"

# --- Expand the AutoFilter range to cover the new rows
$ws.Range("A1:F30").AutoFilter(1)

# --- Update view: freeze pane + selection on OpsTracker
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("E12").Select()

# --- Make OpsTracker the active/selected sheet in the workbook
$wb.Worksheets.Item("InternalAdmin").Select()
$ws.Select()

$wb.Save()
